# "added needed info to RD1"
# Row 6 of Sheet1 is the "Regular US Data" (RD1) test case row; it was
# missing its "input: Time btw births (s)" value. Fill it in.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 3

# Leave the cursor where the editing session ended up.
[void]$ws.Range("C9").Select()
